$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host ($ws.Cells.Item(5,4) | Get-Member | Out-String)
